$d = $word.ActiveDocument

# 1) Merge the two runs about "Para a implementação..." / "Novas estruturas..."
#    into a single run (same combined text, dropping the now-unneeded
#    xml:space="preserve" split point).
$d.Content.Find.Execute(
    "Para a implementação do trabalho decidimos fazer uma modificação profunda no código do trabalho anterior. Novas estruturas de dados necessitavam ser criadas, além de modificações nas estruturas já existentes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Para a implementação do trabalho decidimos fazer uma modificação profunda no código do trabalho anterior. Novas estruturas de dados necessitavam ser criadas, além de modificações nas estruturas já existentes.",
    2
)

# 2) Prepend the new "disco" paragraph content in front of the existing
#    "E por fim, ..." sentence, merging all the small runs that made up
#    that tail into one.
$d.Content.Find.Execute(
    "E por fim, novas funções de impressão (para imprimir o vetor de memória) e calculo de metricas do uso de memória foram implementadas. As figuras abaixo ilustram detalhes citados na implementação:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Além disso, uma implementação que busca simular um disco (assim como recomendado pelo monitor no fórum de dúvidas) foi implementado e sempre que um processo se encontra  bloqueado, tem seu conteudo da memória movido para o vetor de memória em disco. E por fim, novas funções de impressão (para imprimir o vetor de memória) e calculo de metricas do uso de memória foram implementadas. As figuras abaixo ilustram detalhes citados na implementação:",
    2
)
